# ver 1.1 2020/08/09 [430424/32244]
# added BH1750FVI (lux) sensor - prepared (but not moved yet!) new HTML templates
#
# 1) Insert a brand-new worksheet "wire (test)" as the FIRST sheet in the
#    workbook (it becomes the active/selected sheet), with a small table
#    describing the BH1750/MLX test wiring (SDA/SCL/GND/3.3V) and a
#    "Test wire diagram" heading.
# 2) On the existing "list" sheet, document the new sensor's i2c bus /
#    voltage range in the "Спец.блок" row (B27/B28).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "wire (test)" sheet, inserted at the front of the workbook.
# ---------------------------------------------------------------------
$wireWs = $wb.Worksheets.Add()
$wireWs.Name = "wire (test)"

$wireWs.Range("A1").Value = "Test wire diagram"
$wireWs.Range("A1").Font.Bold = $true

$wireWs.Range("A4").Value = "MLX"
$wireWs.Range("A4").Font.Bold = $true
$wireWs.Range("B4").Value = "GND"
$wireWs.Range("E4").Value = "SCL"

$wireWs.Range("B9").Value = "3.3V"
$wireWs.Range("E9").Value = "SDA"

# Small schematic: two "chips" (ovals) joined by wires, echoing the
# hand-drawn wiring diagram from the authored workbook.
$oval1 = $wireWs.Shapes.AddShape(9, 106.96, 40.5, 85.04, 85.04)
$oval1.Name = "Oval 2"
$oval1.Fill.ForeColor.RGB = 4763817

$rect1 = $wireWs.Shapes.AddShape(1, 145, 35.5, 8, 7.5)
$rect1.Name = "Rectangle 3"
$rect1.Fill.ForeColor.RGB = 4763817

$oval2 = $wireWs.Shapes.AddShape(9, 127.5, 59.5, 10.5, 10.5)
$oval2.Name = "Oval 4"
$oval2.Fill.ForeColor.RGB = 14423100

$oval3 = $wireWs.Shapes.AddShape(9, 158.5, 59.5, 10.5, 10.5)
$oval3.Name = "Oval 5"
$oval3.Fill.ForeColor.RGB = 16777215

$oval4 = $wireWs.Shapes.AddShape(9, 127.5, 93.5, 10.5, 10.5)
$oval4.Name = "Oval 6"
$oval4.Fill.ForeColor.RGB = 16777215

$oval5 = $wireWs.Shapes.AddShape(9, 158.5, 93.5, 10.5, 10.5)
$oval5.Name = "Oval 7"
$oval5.Fill.ForeColor.RGB = 16777215

$line1 = $wireWs.Shapes.AddLine(184, 29, 336, 29)
$line1.Name = "Straight Connector 9"
$line1.Line.ForeColor.RGB = 5164313

$line2 = $wireWs.Shapes.AddLine(167, 145, 336, 145)
$line2.Name = "Straight Connector 10"
$line2.Line.ForeColor.RGB = 255

$line3 = $wireWs.Shapes.AddLine(135, 102.5, 167, 145)
$line3.Name = "Straight Connector 11"
$line3.Line.ForeColor.RGB = 255

$line4 = $wireWs.Shapes.AddLine(137, 29, 185, 58)
$line4.Name = "Straight Connector 15"
$line4.Line.ForeColor.RGB = 5164313

$line5 = $wireWs.Shapes.AddLine(170, 58, 339, 58)
$line5.Name = "Straight Connector 18"
$line5.Line.ForeColor.RGB = 13938487

$line6 = $wireWs.Shapes.AddLine(171, 101.5, 340, 101.5)
$line6.Name = "Straight Connector 19"
$line6.Line.ForeColor.RGB = 14423100

# ---------------------------------------------------------------------
# 2. "list" sheet: add the i2c / voltage-range notes for the BH1750.
# ---------------------------------------------------------------------
$listWs = $wb.Worksheets.Item("list")
$listWs.Range("B27").Value = "i2c"
$listWs.Range("B28").Value = "3.3v - 5v"

# ---------------------------------------------------------------------
# 3. Selections: "wire (test)" ends up the active sheet/cell, "list"
#    keeps a B27:B28 selection for the newly documented rows.
# ---------------------------------------------------------------------
$listWs.Range("B27:B28").Select()

$wireWs.Activate()
$wireWs.Range("A13").Select()
